$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-23 Monday" "2026-02-24 Tuesday"

Replace-Text "469÷2=" "492÷9="
Replace-Text "865÷9=" "493÷7="
Replace-Text "813÷6=" "791÷9="
Replace-Text "282÷3=" "815÷8="
Replace-Text "980÷3=" "230÷6="
Replace-Text "887÷3=" "463÷9="
Replace-Text "804÷3=" "729÷8="
Replace-Text "160÷9=" "966÷4="
Replace-Text "669÷6=" "120÷4="
Replace-Text "290÷6=" "383÷9="
Replace-Text "716÷9=" "582÷9="
Replace-Text "694÷2=" "494÷9="
Replace-Text "958÷7=" "266÷6="
Replace-Text "587÷3=" "718÷9="
Replace-Text "535÷4=" "759÷3="
Replace-Text "298÷4=" "562÷3="
Replace-Text "126÷8=" "977÷7="
Replace-Text "519÷6=" "238÷2="
Replace-Text "262÷9=" "769÷4="
Replace-Text "209÷3=" "865÷3="
Replace-Text "606÷8=" "275÷4="
Replace-Text "155÷7=" "650÷2="
Replace-Text "298÷9=" "572÷6="
Replace-Text "440÷7=" "362÷8="
Replace-Text "931÷6=" "173÷7="
